$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows (96, 97) to the "Condicion_Pacientes" table with the
# latest data (1/6/2020 and 2/6/2020), which also grows the table range
# from A1:I95 to A1:I97.

$lo = $ws.ListObjects.Item("Condicion_Pacientes")

$newRow1 = $lo.ListRows.Add()
$newRow2 = $lo.ListRows.Add()

$ws.Range("A96").Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
$ws.Range("B96").Value = 1
$ws.Range("C96").Value = 6
$ws.Range("D96").Value = 2020
$ws.Range("E96").Value = 2374
$ws.Range("F96").Value = 915

$ws.Range("A97").Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
$ws.Range("B97").Value = 2
$ws.Range("C97").Value = 6
$ws.Range("D97").Value = 2020
$ws.Range("E97").Value = 2186
$ws.Range("F97").Value = 1061

# Match the formatting already used throughout the table (column A =
# centered with grey fill matching "Fecha"; B/C/D = centered) by copying
# the style straight from the row above, which already has it applied.
$ws.Range("A95:D95").Copy()
$ws.Range("A96:D97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G97").Select()
